$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A15").Value = 11
$ws.Range("B15").Value = "ค่าวัสดุการศึกษา"
$ws.Range("C15").Value = 1234
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2023-12-30"
